$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 99; this shifts existing rows 99-215 down to 100-216
# and automatically extends the used range / dimension to A1:T216.
$ws.Rows.Item(99).Insert()

# Populate the newly inserted row 99 with the new data record.
$ws.Cells.Item(99, 1).Value = 7
$ws.Cells.Item(99, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(99, 3).Value = "Ñuble"
$ws.Cells.Item(99, 4).Value = 44664
$ws.Cells.Item(99, 5).Value = 16
$ws.Cells.Item(99, 6).Value = "Fruta"
$ws.Cells.Item(99, 7).Value = 100108
$ws.Cells.Item(99, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(99, 9).Value = 100108005
$ws.Cells.Item(99, 10).Value = "Piña"
$ws.Cells.Item(99, 11).Value = "Caramelo"
$ws.Cells.Item(99, 12).Value = "Segunda"
$ws.Cells.Item(99, 13).Value = 120
$ws.Cells.Item(99, 14).Value = 15000
$ws.Cells.Item(99, 15).Value = 16000
$ws.Cells.Item(99, 16).Value = 15500
$ws.Cells.Item(99, 17).Value = "$/caja 14 unidades"
$ws.Cells.Item(99, 18).Value = "Ecuador"
$ws.Cells.Item(99, 19).Value = 1107
$ws.Cells.Item(99, 20).Value = 14
